$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 50518.99
$ws.Range("D2").Value = 141789.57
$ws.Range("E2").Value = 162819.35
$ws.Range("F2").Value = 331129.46
$ws.Range("G2").Value = 402689.42
$ws.Range("C3").Value = 20354.44
$ws.Range("D3").Value = 150139.79
$ws.Range("E3").Value = 187082.86
$ws.Range("F3").Value = 343644.5
$ws.Range("G3").Value = 455942.8
$ws.Range("D4").Value = 166638.29
$ws.Range("E4").Value = 193765.12
$ws.Range("F4").Value = 323694.55
$ws.Range("G4").Value = 380953.36
$ws.Range("D5").Value = 130298.05
$ws.Range("E5").Value = 215432.77
$ws.Range("F5").Value = 396044.77
$ws.Range("G5").Value = 424509.78
$ws.Range("D6").Value = 115915.98
$ws.Range("E6").Value = 220684.43
$ws.Range("F6").Value = 378096.36
$ws.Range("G6").Value = 490485.37
$ws.Range("D7").Value = 141564.4
$ws.Range("E7").Value = 246532.18
$ws.Range("F7").Value = 403206.55
$ws.Range("G7").Value = 126538.86
$ws.Range("D8").Value = 136354.87
$ws.Range("E8").Value = 224304.98
$ws.Range("F8").Value = 389318.8
$ws.Range("C9").Value = 117632.99
$ws.Range("D9").Value = 148788.4
$ws.Range("E9").Value = 267309.71
$ws.Range("F9").Value = 381906.92
$ws.Range("B10").Value = 40488.91
$ws.Range("D10").Value = 153181.38
$ws.Range("E10").Value = 292606.09
$ws.Range("F10").Value = 386980.51
$ws.Range("B11").Value = 33106.69
$ws.Range("D11").Value = 134868.15
$ws.Range("E11").Value = 331822.65
$ws.Range("F11").Value = 418014.5
$ws.Range("D12").Value = 151292.14
$ws.Range("E12").Value = 249653.56
$ws.Range("F12").Value = 309538.28
$ws.Range("B13").Value = 23097.81
$ws.Range("D13").Value = 201369.38
$ws.Range("E13").Value = 251128.74
$ws.Range("F13").Value = 381502.92
